$wb = $excel.ActiveWorkbook

# Row -> new value for column F ("想去人数") that changed in this update.
$changes = @{
    2  = 2081
    4  = 131
    5  = 44
    7  = 1723
    9  = 709
    14 = 234
    16 = 142
    17 = 117
    19 = 3994
    22 = 448
    23 = 383
    24 = 944
    25 = 761
    27 = 13
    28 = 36
    29 = 1793
    30 = 34
    32 = 69
}

# Both "展览" and "全部类型" sheets carry the same table and both were
# updated identically in the source data refresh.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $changes.Keys) {
        $ws.Cells.Item($row, 6).Value = $changes[$row]
    }
}
